$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.029585357532213
    "D2" = 1.031965529203292
    "E2" = 1.02934604942493
    "F2" = 1.038377287140346
    "I2" = 1.031754744203483
    "J2" = 1.034731317919554
    "K2" = 1.034772123665848
    "L2" = 1.032160227836086
    "M2" = 1.041165493558144
    "B3" = 1.02
    "C3" = 1.0310575500719
    "D3" = 1.033028217522837
    "E3" = 1.030613589053388
    "F3" = 1.039955187531475
    "I3" = 1.032068959609567
    "J3" = 1.035841933759487
    "K3" = 1.035642805629568
    "L3" = 1.033234656961156
    "M3" = 1.042551367157641
    "B4" = 1.02
    "C4" = 1.0320086449775
    "D4" = 1.033714378219901
    "E4" = 1.03143268954541
    "F4" = 1.040974837827389
    "I4" = 1.032270253801178
    "J4" = 1.036558705696557
    "K4" = 1.036204136795907
    "L4" = 1.033928259641451
    "M4" = 1.043446286430386
    "B5" = 1.02
    "C5" = 1.03240813086359
    "D5" = 1.034002492719051
    "E5" = 1.031776785811581
    "F5" = 1.041403182251292
    "I5" = 1.032354395086768
    "J5" = 1.036859594230037
    "K5" = 1.036439631829919
    "L5" = 1.034219466319082
    "M5" = 1.043822078907443
    "B6" = 1.02
    "C6" = 1.032475185738517
    "D6" = 1.034050848163149
    "E6" = 1.031834546379264
    "F6" = 1.04147508484053
    "I6" = 1.032368494519753
    "J6" = 1.036910088955411
    "K6" = 1.036479143951878
    "L6" = 1.034268338837552
    "M6" = 1.043885151023347
    "B7" = 1.02
    "C7" = 1.032013984310681
    "D7" = 1.033718229379803
    "E7" = 1.031437288366107
    "F7" = 1.040980562620215
    "I7" = 1.032271379995987
    "J7" = 1.0365627279154
    "K7" = 1.036207285406767
    "L7" = 1.033932152261473
    "M7" = 1.043451309474577
    "B8" = 1.02
    "C8" = 1.030083209567404
    "D8" = 1.032324976662308
    "E8" = 1.029774647846904
    "F8" = 1.038910831874038
    "I8" = 1.031861354981383
    "J8" = 1.035107046407274
    "K8" = 1.035066803412009
    "L8" = 1.032523675692274
    "M8" = 1.041634240150849
    "B9" = 1.02
    "C9" = 1.026669015588568
    "D9" = 1.029858433108188
    "E9" = 1.026836323082767
    "F9" = 1.035252939656885
    "I9" = 1.031123254302931
    "J9" = 1.032527375057895
    "K9" = 1.033041176492793
    "L9" = 1.030029086613864
    "M9" = 1.038417956368884
    "B10" = 1.02
    "C10" = 1.02438439114048
    "D10" = 1.028206089816871
    "E10" = 1.024871346474502
    "F10" = 1.032806588049648
    "I10" = 1.030620594051656
    "J10" = 1.030797447820595
    "K10" = 1.031679768205416
    "L10" = 1.028357183384409
    "M10" = 1.036263635243923
    "B11" = 1.02
    "C11" = 1.023393001303178
    "D11" = 1.027488650881246
    "E11" = 1.024018965605556
    "F11" = 1.03174533444322
    "I11" = 1.030400397988559
    "J11" = 1.030045885674649
    "K11" = 1.031087596895345
    "L11" = 1.027631060628078
    "M11" = 1.035328282748688
    "B12" = 1.02
    "C12" = 1.023024425284482
    "D12" = 1.02722186160669
    "E12" = 1.023702116182853
    "F12" = 1.031350831929353
    "I12" = 1.030318223319004
    "J12" = 1.029766340790519
    "K12" = 1.030867231346013
    "L12" = 1.027361012803858
    "M12" = 1.034980463910072
    "B13" = 1.02
    "C13" = 1.023103501161422
    "D13" = 1.027279102499817
    "E13" = 1.023770092305442
    "F13" = 1.031435468053407
    "I13" = 1.03033586748089
    "J13" = 1.029826321478956
    "K13" = 1.030914518995698
    "L13" = 1.027418954172862
    "M13" = 1.035055089903407
    "B14" = 1.02
    "C14" = 1.023362541474316
    "D14" = 1.027466604159089
    "E14" = 1.023992779612659
    "F14" = 1.031712731036809
    "I14" = 1.030393613252647
    "J14" = 1.030022786224098
    "K14" = 1.031069389738522
    "L14" = 1.027608745230884
    "M14" = 1.03529953988813
    "B15" = 1.02
    "C15" = 1.023522100792903
    "D15" = 1.02758209013281
    "E15" = 1.024129952971578
    "F15" = 1.031883521103053
    "I15" = 1.030429141368905
    "J15" = 1.030143783915054
    "K15" = 1.031164756657562
    "L15" = 1.027725637382841
    "M15" = 1.035450102000655
    "B16" = 1.02
    "C16" = 1.02445014069788
    "D16" = 1.028253662038787
    "E16" = 1.024927883299134
    "F16" = 1.03287697764741
    "I16" = 1.030635154009346
    "J16" = 1.030847273436553
    "K16" = 1.03171901191154
    "L16" = 1.028405327340489
    "M16" = 1.036325657699176
    "B17" = 1.02
    "C17" = 1.0250316981533
    "D17" = 1.028674391807863
    "E17" = 1.025427988626992
    "F17" = 1.033499612891701
    "I17" = 1.030763698410493
    "J17" = 1.031287881795539
    "K17" = 1.032065962390351
    "L17" = 1.028831091423682
    "M17" = 1.0368741903283
    "B18" = 1.02
    "C18" = 1.025370705668897
    "D18" = 1.028919607302782
    "E18" = 1.025719544157829
    "F18" = 1.033862596375861
    "I18" = 1.030838431186671
    "J18" = 1.031544641195683
    "K18" = 1.032268075230538
    "L18" = 1.029079222732322
    "M18" = 1.037193897872349
    "B19" = 1.02
    "C19" = 1.025486263961713
    "D19" = 1.029003187547872
    "E19" = 1.025818932257436
    "F19" = 1.03398633259034
    "I19" = 1.03086387164562
    "J19" = 1.031632149012392
    "K19" = 1.032336946966356
    "L19" = 1.02916379372198
    "M19" = 1.037302869065345
    "B20" = 1.02
    "C20" = 1.024969323825083
    "D20" = 1.028629271080999
    "E20" = 1.025374347359255
    "F20" = 1.033432829657383
    "I20" = 1.030749932170865
    "J20" = 1.031240633562287
    "K20" = 1.032028764560327
    "L20" = 1.028785432700837
    "M20" = 1.036815363092185
    "B21" = 1.02
    "C21" = 1.023286269734947
    "D21" = 1.027411397934027
    "E21" = 1.023927210355002
    "F21" = 1.03163109251852
    "I21" = 1.03037661918241
    "J21" = 1.029964942832469
    "K21" = 1.031023795435649
    "L21" = 1.027552865752172
    "M21" = 1.035227566224446
    "B22" = 1.02
    "C22" = 1.022226154661402
    "D22" = 1.02664393131335
    "E22" = 1.023015961968546
    "F22" = 1.030496495272163
    "I22" = 1.030139679348135
    "J22" = 1.029160656287495
    "K22" = 1.030389575451595
    "L22" = 1.026775969272901
    "M22" = 1.034227010102582
    "B23" = 1.02
    "C23" = 1.022788326103571
    "D23" = 1.027050946818111
    "E23" = 1.023499164717172
    "F23" = 1.031098138255298
    "I23" = 1.030265497142941
    "J23" = 1.029587235492507
    "K23" = 1.030726012548531
    "L23" = 1.027188002124972
    "M23" = 1.034757639707309
    "B24" = 1.02
    "C24" = 1.024997508733108
    "D24" = 1.028649659779711
    "E24" = 1.025398585995032
    "F24" = 1.033463006710883
    "I24" = 1.030756153299114
    "J24" = 1.031261983746907
    "K24" = 1.032045573454416
    "L24" = 1.02880606456043
    "M24" = 1.036841945336743
    "B25" = 1.02
    "C25" = 1.027553128538658
    "D25" = 1.030497479767422
    "E25" = 1.027596998509573
    "F25" = 1.036199921239623
    "I25" = 1.031315929619931
    "J25" = 1.033196044668069
    "K25" = 1.033566767144221
    "L25" = 1.039251194214634
    "M25" = 1.043142645695669
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"